$d = $word.ActiveDocument

# Locate the three paragraphs to remove:
#   1. the (empty) paragraph right after "Consolidação das Leis do Trabalho..."
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages.
#      Original theme under Creative Commons Attribution"
# and delete that whole span (including their paragraph marks), leaving the
# bibliography entry paragraph and the trailing blank paragraph untouched.

$copyrightMark = [char]0xA9
$startMarker = "Ver no Jupiter Salvar em pdf Salvar em docx"
$endMarker = $copyrightMark + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$count = $d.Paragraphs.Count
$startIndex = -1
$endIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "*$startMarker*") {
        $startIndex = $i
    }
    if ($text -like "*$endMarker*") {
        $endIndex = $i
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    # The blank paragraph immediately preceding "Ver no Jupiter..." is also removed.
    $deleteFrom = $startIndex - 1

    $rangeStart = $d.Paragraphs.Item($deleteFrom).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
